# ESPNSportsAnalytics.xlsx - "Updated for 2018 INFORMS Annual Meeting"
#
# The "Main" worksheet holds year-by-year playoff (PO) / home-field (HF)
# columns for every team, grouped by sport:
#   E/F = 2014PO/2014HF   G/H = 2015PO/2015HF   I/J = 2016PO/2016HF
#   K/L = 2017PO/2017HF   M/N = 2018PO/2018HF
#
# For this refresh:
#   * MLB rows (2018 season just completed) get their previously-empty
#     2018PO/2018HF (M/N) columns filled in with the real results.
#   * NBA rows still had their data living one year "behind" (their 2014
#     columns were never populated), so the whole block of year columns
#     shifts left by one year-pair (2015->2014, 2016->2015, 2017->2016,
#     2018->2017) and the now-unused 2018 columns (M/N) are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# MLB: populate the new 2018PO / 2018HF (M/N) results
# ---------------------------------------------------------------------
$mlb2018 = @{
    2  = @(1, 5)
    4  = @(1, 2)
    5  = @(1, 3)
    6  = @(1, 2)
    7  = @(1, 1)
    13 = @(1, 4)
    20 = @(1, 3)
    25 = @(1, 2)
}

foreach ($row in $mlb2018.Keys) {
    $vals = $mlb2018[$row]
    $ws.Range("M$row").Value = $vals[0]
    $ws.Range("N$row").Value = $vals[1]
}

# ---------------------------------------------------------------------
# NBA: shift the year-pair columns left by one year, dropping the old
# 2014 data and clearing the trailing 2018 columns.
# ---------------------------------------------------------------------
$yearCols = @("E", "F"), @("G", "H"), @("I", "J"), @("K", "L"), @("M", "N")

for ($row = 32; $row -le 58; $row++) {
    $original = @()
    foreach ($pair in $yearCols) {
        $poVal = $ws.Range("$($pair[0])$row").Value2
        $hfVal = $ws.Range("$($pair[1])$row").Value2
        $original += , @($poVal, $hfVal)
    }

    for ($i = 0; $i -lt $yearCols.Length; $i++) {
        $pair = $yearCols[$i]
        if ($i -lt $yearCols.Length - 1) {
            $newVal = $original[$i + 1]
        } else {
            $newVal = @(0, 0)
        }
        $ws.Range("$($pair[0])$row").Value = $newVal[0]
        $ws.Range("$($pair[1])$row").Value = $newVal[1]
    }
}

# ---------------------------------------------------------------------
# Update the active cell / selection left over from editing
# ---------------------------------------------------------------------
$ws.Range("N8").Select()
